# Upload new version with timestamp
# A new product row ("LYSE 0.65% NASAL DROPS 15 ML") is inserted into the
# sorted product list at row 33 (alphabetically between "LEZBERG TRIO..."
# and "MAXILASE..."). All existing rows from 33..61 shift their content
# down by one row (34..62); the "total" row and the footer row move from
# 62/63 to 63/64; the running total is increased by the new row's price.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: push the two trailing "summary" rows further down by one ---
# Old row 63 (footer, e.g. timestamp / page) -> new row 64
$ws.Range("A63:N63").Copy()
$ws.Range("A64:N64").PasteSpecial(-4122)
$ws.Range("A63:N63").Copy()
$ws.Range("A64:N64").PasteSpecial(-4163)
$ws.Range("A64:E64").Merge()
$ws.Range("F64:G64").Merge()
$ws.Range("I64:N64").Merge()
$ws.Rows.Item(64).RowHeight = 17.25

# Old row 62 (running total) -> new row 63
$ws.Range("A62:N62").Copy()
$ws.Range("A63:N63").PasteSpecial(-4122)
$ws.Range("A62:N62").Copy()
$ws.Range("A63:N63").PasteSpecial(-4163)
$ws.Range("K63:N63").Merge()
$ws.Rows.Item(63).RowHeight = 25.5

# --- Step 2: shift data rows 33..61 down to 34..62 (iterate from the
#     bottom up so we never overwrite a row before it has been copied).
#     Column A is the row's fixed sequence number (row-3) and must stay
#     put, so only B..N are copied. ---
for ($r = 61; $r -ge 33; $r--) {
    $dst = $r + 1
    $ws.Range("B$r`:N$r").Copy()
    $ws.Range("B$dst`:N$dst").PasteSpecial(-4122)
    $ws.Range("B$r`:N$r").Copy()
    $ws.Range("B$dst`:N$dst").PasteSpecial(-4163)
    $ws.Range("B$dst`:G$dst").Merge()
    $ws.Range("H$dst`:K$dst").Merge()
    $ws.Range("L$dst`:M$dst").Merge()
}
$ws.Range("A62").Value = 59

# --- Step 3: write the brand-new product row at row 33 ---
$ws.Range("B33").Value = "LYSE 0.65% NASAL DROPS 15 ML"
$ws.Range("H33").Value = "6:0"
$ws.Range("L33").Value = 15
$ws.Range("N33").Value = "1:0"

# --- Step 4: update the running total (old total + new row's price) ---
$ws.Range("K63").Value = 3257.6399999999999
